$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold a given string as TEXT (avoids Excel auto-converting
# numeric-looking strings like "173.65" into real numbers), then strip the
# temporary Text number-format back off so no stray style index is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '64.321.31'
$ws.Range('E2').Value = '  +0.37%  '

Set-TextValue $ws.Range('D3') '3.330.33'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('E4').Value = '  +0.13%  '

Set-TextValue $ws.Range('D5') '552.87'
$ws.Range('E5').Value = '  +0.38%  '

Set-TextValue $ws.Range('D6') '173.65'
$ws.Range('E6').Value = '  +0.49%  '

Set-TextValue $ws.Range('D7') '0.630'
$ws.Range('E7').Value = '  +2.74%  '

$ws.Range('E8').Value = '  -0.02%  '

Set-TextValue $ws.Range('D9') '3.322.39'
$ws.Range('E9').Value = '  +0.17%  '

Set-TextValue $ws.Range('D10') '0.171'
$ws.Range('E10').Value = '  +6.22%  '

Set-TextValue $ws.Range('D11') '0.633'
$ws.Range('E11').Value = '  +2.06%  '

Set-TextValue $ws.Range('D12') '53.48'
$ws.Range('E12').Value = '  +0.66%  '

Set-TextValue $ws.Range('D13') '0.0000278'
$ws.Range('E13').Value = '  +2.27%  '

Set-TextValue $ws.Range('D14') '9.10'
$ws.Range('E14').Value = '  +1.16%  '

Set-TextValue $ws.Range('D15') '3.863.06'
$ws.Range('E15').Value = '  +0.12%  '

$ws.Range('E16').Value = '  +3.08%  '

Set-TextValue $ws.Range('D17') '18.15'
$ws.Range('E17').Value = '  -0.56%  '

Set-TextValue $ws.Range('D18') '3.327.94'
$ws.Range('E18').Value = '  +0.10%  '

Set-TextValue $ws.Range('D19') '64.432.48'
$ws.Range('E19').Value = '  +0.83%  '

Set-TextValue $ws.Range('D20') '11.76'
$ws.Range('E20').Value = '  +0.32%  '

$ws.Range('E21').Value = '  +1.69%  '

Set-TextValue $ws.Range('D22') '452.52'
$ws.Range('E22').Value = '  +6.37%  '

Set-TextValue $ws.Range('D23') '5.05'
$ws.Range('E23').Value = '  +9.43%  '

$ws.Range('E24').Value = '  -0.26%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D25') '87.69'
$ws.Range('E25').Value = '  +4.40%  '

$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D26') '13.90'
$ws.Range('E26').Value = '  +4.78%  '

$ws.Range('E27').Value = '  +2.94%  '

Set-TextValue $ws.Range('D28') '10.59'
$ws.Range('E28').Value = '  -0.46%  '

Set-TextValue $ws.Range('D29') '31.15'
$ws.Range('E29').Value = '  +5.03%  '

$ws.Range('E30').Value = '  +0.11%  '

Set-TextValue $ws.Range('D31') '6.52'
$ws.Range('E31').Value = '  -2.14%  '

$ws.Range('E32').Value = '  +0.37%  '

$ws.Range('E33').Value = '  +6.11%  '

Set-TextValue $ws.Range('D34') '567.90'
$ws.Range('E34').Value = '  -4.36%  '

$ws.Range('E35').Value = '  +0.54%  '

$ws.Range('E37').Value = '  -0.91%  '

Set-TextValue $ws.Range('D38') '3.51'
$ws.Range('E38').Value = '  +0.51%  '

Set-TextValue $ws.Range('D39') '35.39'
$ws.Range('E39').Value = '  +0.20%  '

Set-TextValue $ws.Range('D40') '0.367'
$ws.Range('E40').Value = '  +0.70%  '

Set-TextValue $ws.Range('D41') '0.0₃0730'
$ws.Range('E41').Value = '  -2.42%  '

Set-TextValue $ws.Range('D42') '3.066.18'
$ws.Range('E42').Value = '  -0.79%  '

Set-TextValue $ws.Range('D43') '0.0415'
$ws.Range('E43').Value = '  +2.76%  '

$ws.Range('E44').Value = '  -1.01%  '

Set-TextValue $ws.Range('D45') '2.45'
$ws.Range('E45').Value = '  +0.55%  '

$ws.Range('E46').Value = '  +4.23%  '

Set-TextValue $ws.Range('D47') '3.18'
$ws.Range('E47').Value = '  -1.45%  '

Set-TextValue $ws.Range('D48') '1.00'
$ws.Range('E48').Value = '  +0.19%  '

Set-TextValue $ws.Range('D49') '140.97'
$ws.Range('E49').Value = '  +6.33%  '

Set-TextValue $ws.Range('D50') '2.51'
$ws.Range('E50').Value = '  -2.46%  '

$ws.Range('E51').Value = '  +0.33%  '
